# Insert a new record at row 80 (Macroferia Regional de Talca - Ají),
# shifting all existing rows 80..171 down by one (to 81..172).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("80").Insert()

$ws.Range("A80").Value = 5
$ws.Range("B80").Value = "Macroferia Regional de Talca"
$ws.Range("C80").Value = "Maule"
$ws.Range("D80").Value = 44586
$ws.Range("E80").Value = 7
$ws.Range("F80").Value = 100112021
$ws.Range("G80").Value = "Ají"
$ws.Range("H80").Value = "Americana (o)"
$ws.Range("I80").Value = "Primera"
$ws.Range("J80").Value = 150
$ws.Range("K80").Value = 20000
$ws.Range("L80").Value = 20000
$ws.Range("M80").Value = 20000
$ws.Range("N80").Value = "`$/saco 25 kilos"
$ws.Range("O80").Value = "Región del Maule"
$ws.Range("P80").Value = 800
$ws.Range("Q80").Value = 25
$ws.Range("R80").Value = "Hortaliza"
